$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove "Favorable rate" row (Recession rate criteria section) ---
$ws.Rows.Item(15).Delete()

# --- Remove the whole "Elevation criteria" section (header + 3 rows) ---
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(17).Delete()

# --- Remove "Favorable inundation" row ---
$ws.Rows.Item(18).Delete()

# --- Update the remaining "Recession rate criteria" values (now rows 15-16) ---
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 5
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 10

# --- Update the "Bed shear stress criteria" D-column values ---
$ws.Range("D11").Value = 0.01
$ws.Range("D12").Value = 0.005

# --- Rename + restyle the "Inundation criteria" header row (now row 17) ---
$ws.Range("A17").Value = "Inundation criteria (partial or complete shoot inundation)"
$a17 = $ws.Range("A17")
$a17.WrapText = $true
$ws.Rows.Item(17).RowHeight = 43.2

# --- Update the remaining inundation rows (now rows 18-19) ---
$ws.Range("C18").Value = 14
$ws.Range("D18").Value = 14
$ws.Range("C19").Value = 28
$ws.Range("D19").Value = 28

# --- Restore the selection to match the saved view ---
$ws.Range("E13").Select() | Out-Null
